$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; this shifts the existing rows 42..154
# down to 43..155 (matching the new dimension A1:R155).
$ws.Rows("42:42").Insert()

# Populate the newly inserted row 42 with its data.
$ws.Range("A42").Value2 = 10
$ws.Range("B42").Value2 = "Vega Modelo de Temuco"
$ws.Range("C42").Value2 = "La Araucanía"
$ws.Range("D42").Value2 = 44925
$ws.Range("E42").Value2 = 9
$ws.Range("F42").Value2 = 100114002
$ws.Range("G42").Value2 = "Camote"
$ws.Range("H42").Value2 = "Sin especificar"
$ws.Range("I42").Value2 = "Primera"
$ws.Range("J42").Value2 = 10
$ws.Range("K42").Value2 = 24000
$ws.Range("L42").Value2 = 24000
$ws.Range("M42").Value2 = 24000
$ws.Range("N42").Value2 = "$/malla 20 kilos"
$ws.Range("O42").Value2 = "Perú"
$ws.Range("P42").Value2 = 1200
$ws.Range("Q42").Value2 = 20
$ws.Range("R42").Value2 = "Hortaliza"

# Match the date-formatted style used by the other rows in column D.
$ws.Range("D42").NumberFormat = $ws.Range("D43").NumberFormat
